# Leave Balance - Hourly Run
# Updates the "RunMode" column (C) on the All_Scenarios sheet:
#   - Rows 2-90   -> "No"
#   - Rows 91-125 -> "Yes"
# and moves the active selection down to the newly "Yes" block (C91:C125).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All_Scenarios")

# Rows 2 through 90 become "No"
$ws.Range("C2:C90").Value = "No"

# Rows 91 through 125 become "Yes"
$ws.Range("C91:C125").Value = "Yes"

# Move the visible selection to match the new block of "Yes" rows
[void]$ws.Range("C91:C125").Select()
